# Auto-generated PowerShell-style Excel COM-interop script
# Applies the cell-value updates from the "Updated symbol list" commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D (Price) and E (Volume(1h)) hold numeric-looking text
# (e.g. "318.57", "4.71%"). Excel auto-converts such strings typed into a
# ".Value" assignment into real numbers/percentages, which would change the
# cell type away from the plain text the workbook uses. Setting the number
# format to Text ("@") before assigning the value keeps these as literal text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "318.57"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.71%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.14"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.41%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.175"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.37%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08231"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "5.08%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.151"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.008"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.55%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9275"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.84%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1019"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1891"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.53%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09204"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "6.33%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03622"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.15%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09924"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.28%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001446"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.54%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005709"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.22%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.463"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.02%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.59%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.801"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "17.93%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.52%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.214"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "7.76%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1301"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.07%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2192"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.64%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04597"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.57%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001248"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004731"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.73%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-11.46%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004505"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-5.39%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02010"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "9.35%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04929"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.38%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007768"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.52%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1401"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.03%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007824"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.02%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002098"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.15%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01194"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "8.56%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006460"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.53%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.16%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "31.16"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-34.94%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001902"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-5.12%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.16%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.16%"
